$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 obsolete data rows (previously rows 17-21, MuSCs-target rows trimmed
# and Resolving-Mac-source block collapsed) so the table now spans A1:T16
$ws.Range("A17:T21").EntireRow.Delete()

# Refresh remaining data rows (2-16) with updated TPM-based statistics
# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "Itgax"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.820425
$ws.Range("H2").Value = 11.461275
$ws.Range("I2").Value = 0.02049663039797357
$ws.Range("J2").Value = 0.02049663039797357
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02697933333333333
$ws.Range("N2").Value = 0.080938
$ws.Range("O2").Value = 0.0003365168416393062
$ws.Range("P2").Value = 0.0003365168416393062
$ws.Range("Q2").Value = 0.10307251955
$ws.Range("R2").Value = 0.92765267595
$ws.Range("S2").Value = 0.00000689746132577426
$ws.Range("T2").Value = 0.000006897461325774261

# Row 3: ECs -> Inflammatory-Mac
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "Itgax"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.820425
$ws.Range("H3").Value = 11.461275
$ws.Range("I3").Value = 0.02049663039797357
$ws.Range("J3").Value = 0.02049663039797357
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 40.972402
$ws.Range("N3").Value = 122.917206
$ws.Range("O3").Value = 0.5110542630933305
$ws.Range("P3").Value = 0.5110542630933306
$ws.Range("Q3").Value = 156.53198891085
$ws.Range("R3").Value = 1408.78790019765
$ws.Range("S3").Value = 0.01047489034393274
$ws.Range("T3").Value = 0.01047489034393274

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "Itgax"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.820425
$ws.Range("H4").Value = 11.461275
$ws.Range("I4").Value = 0.02049663039797357
$ws.Range("J4").Value = 0.02049663039797357
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 39.172931
$ws.Range("N4").Value = 117.518793
$ws.Range("O4").Value = 0.4886092200650302
$ws.Range("P4").Value = 0.4886092200650302
$ws.Range("Q4").Value = 149.657244915675
$ws.Range("R4").Value = 1346.915204241075
$ws.Range("S4").Value = 0.01001484259271505
$ws.Range("T4").Value = 0.01001484259271505

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "Itgax"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 140.5890806666667
$ws.Range("H5").Value = 421.767242
$ws.Range("I5").Value = 0.7542622677884155
$ws.Range("J5").Value = 0.7542622677884157
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02697933333333333
$ws.Range("N5").Value = 0.080938
$ws.Range("O5").Value = 0.0003365168416393062
$ws.Range("P5").Value = 0.0003365168416393062
$ws.Range("Q5").Value = 3.792999670332888
$ws.Range("R5").Value = 34.136997032996
$ws.Range("S5").Value = 0.0002538219561238582
$ws.Range("T5").Value = 0.0002538219561238582

# Row 6: FAPs -> Inflammatory-Mac
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "Itgax"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 140.5890806666667
$ws.Range("H6").Value = 421.767242
$ws.Range("I6").Value = 0.7542622677884155
$ws.Range("J6").Value = 0.7542622677884157
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 40.972402
$ws.Range("N6").Value = 122.917206
$ws.Range("O6").Value = 0.5110542630933305
$ws.Range("P6").Value = 0.5110542630933306
$ws.Range("Q6").Value = 5760.272329885095
$ws.Range("R6").Value = 51842.45096896586
$ws.Range("S6").Value = 0.385468947443713
$ws.Range("T6").Value = 0.3854689474437131

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "Itgax"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 140.5890806666667
$ws.Range("H7").Value = 421.767242
$ws.Range("I7").Value = 0.7542622677884155
$ws.Range("J7").Value = 0.7542622677884157
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 39.172931
$ws.Range("N7").Value = 117.518793
$ws.Range("O7").Value = 0.4886092200650302
$ws.Range("P7").Value = 0.4886092200650302
$ws.Range("Q7").Value = 5507.286356308767
$ws.Range("R7").Value = 49565.5772067789
$ws.Range("S7").Value = 0.3685394983885786
$ws.Range("T7").Value = 0.3685394983885787

# Row 8: Inflammatory-Mac -> ECs
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "C3"
$ws.Range("C8").Value = "Itgax"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 30.51067
$ws.Range("H8").Value = 91.53201
$ws.Range("I8").Value = 0.1636901460399144
$ws.Range("J8").Value = 0.1636901460399144
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02697933333333333
$ws.Range("N8").Value = 0.080938
$ws.Range("O8").Value = 0.0003365168416393062
$ws.Range("P8").Value = 0.0003365168416393062
$ws.Range("Q8").Value = 0.8231575361533333
$ws.Range("R8").Value = 7.40841782538
$ws.Range("S8").Value = 0.00005508449095282879
$ws.Range("T8").Value = 0.0000550844909528288

# Row 9: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "C3"
$ws.Range("C9").Value = "Itgax"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 30.51067
$ws.Range("H9").Value = 91.53201
$ws.Range("I9").Value = 0.1636901460399144
$ws.Range("J9").Value = 0.1636901460399144
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 40.972402
$ws.Range("N9").Value = 122.917206
$ws.Range("O9").Value = 0.5110542630933305
$ws.Range("P9").Value = 0.5110542630933306
$ws.Range("Q9").Value = 1250.09543652934
$ws.Range("R9").Value = 11250.85892876406
$ws.Range("S9").Value = 0.08365454696006812
$ws.Range("T9").Value = 0.08365454696006813

# Row 10: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "C3"
$ws.Range("C10").Value = "Itgax"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 30.51067
$ws.Range("H10").Value = 91.53201
$ws.Range("I10").Value = 0.1636901460399144
$ws.Range("J10").Value = 0.1636901460399144
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 39.172931
$ws.Range("N10").Value = 117.518793
$ws.Range("O10").Value = 0.4886092200650302
$ws.Range("P10").Value = 0.4886092200650302
$ws.Range("Q10").Value = 1195.19237067377
$ws.Range("R10").Value = 10756.73133606393
$ws.Range("S10").Value = 0.07998051458889348
$ws.Range("T10").Value = 0.07998051458889348

# Row 11: MuSCs -> ECs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "C3"
$ws.Range("C11").Value = "Itgax"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.258813
$ws.Range("H11").Value = 0.776439
$ws.Range("I11").Value = 0.001388535150720334
$ws.Range("J11").Value = 0.001388535150720334
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.02697933333333333
$ws.Range("N11").Value = 0.080938
$ws.Range("O11").Value = 0.0003365168416393062
$ws.Range("P11").Value = 0.0003365168416393062
$ws.Range("Q11").Value = 0.006982602198
$ws.Range("R11").Value = 0.06284341978199999
$ws.Range("S11").Value = 0.0000004672654634255648
$ws.Range("T11").Value = 0.0000004672654634255648

# Row 12: MuSCs -> Inflammatory-Mac
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "C3"
$ws.Range("C12").Value = "Itgax"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.258813
$ws.Range("H12").Value = 0.776439
$ws.Range("I12").Value = 0.001388535150720334
$ws.Range("J12").Value = 0.001388535150720334
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 40.972402
$ws.Range("N12").Value = 122.917206
$ws.Range("O12").Value = 0.5110542630933305
$ws.Range("P12").Value = 0.5110542630933306
$ws.Range("Q12").Value = 10.604190278826
$ws.Range("R12").Value = 95.43771250943401
$ws.Range("S12").Value = 0.0007096168082305669
$ws.Range("T12").Value = 0.0007096168082305669

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "C3"
$ws.Range("C13").Value = "Itgax"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.258813
$ws.Range("H13").Value = 0.776439
$ws.Range("I13").Value = 0.001388535150720334
$ws.Range("J13").Value = 0.001388535150720334
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 39.172931
$ws.Range("N13").Value = 117.518793
$ws.Range("O13").Value = 0.4886092200650302
$ws.Range("P13").Value = 0.4886092200650302
$ws.Range("Q13").Value = 10.138463790903
$ws.Range("R13").Value = 91.24617411812699
$ws.Range("S13").Value = 0.0006784510770263416
$ws.Range("T13").Value = 0.0006784510770263415

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "C3"
$ws.Range("C14").Value = "Itgax"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 11.213844
$ws.Range("H14").Value = 33.641532
$ws.Range("I14").Value = 0.0601624206229761
$ws.Range("J14").Value = 0.0601624206229761
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.02697933333333333
$ws.Range("N14").Value = 0.080938
$ws.Range("O14").Value = 0.0003365168416393062
$ws.Range("P14").Value = 0.0003365168416393062
$ws.Range("Q14").Value = 0.302542035224
$ws.Range("R14").Value = 2.722878317016
$ws.Range("S14").Value = 0.00002024566777341938
$ws.Range("T14").Value = 0.00002024566777341938

# Row 15: Resolving-Mac -> Inflammatory-Mac
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "C3"
$ws.Range("C15").Value = "Itgax"
$ws.Range("D15").Value = "Inflammatory-Mac"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 11.213844
$ws.Range("H15").Value = 33.641532
$ws.Range("I15").Value = 0.0601624206229761
$ws.Range("J15").Value = 0.0601624206229761
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 40.972402
$ws.Range("N15").Value = 122.917206
$ws.Range("O15").Value = 0.5110542630933305
$ws.Range("P15").Value = 0.5110542630933306
$ws.Range("Q15").Value = 459.458124333288
$ws.Range("R15").Value = 4135.123118999592
$ws.Range("S15").Value = 0.03074626153738604
$ws.Range("T15").Value = 0.03074626153738604

# Row 16: Resolving-Mac -> Resolving-Mac
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "C3"
$ws.Range("C16").Value = "Itgax"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 11.213844
$ws.Range("H16").Value = 33.641532
$ws.Range("I16").Value = 0.0601624206229761
$ws.Range("J16").Value = 0.0601624206229761
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 39.172931
$ws.Range("N16").Value = 117.518793
$ws.Range("O16").Value = 0.4886092200650302
$ws.Range("P16").Value = 0.4886092200650302
$ws.Range("Q16").Value = 439.279137256764
$ws.Range("R16").Value = 3953.512235310875
$ws.Range("S16").Value = 0.02939591341781664
$ws.Range("T16").Value = 0.02939591341781663

Write-Output $ws.UsedRange.Address()
